# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") with the latest scraped values. Cells in D that look like
# plain numbers ("1.003", "309.22", ...) are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's original
# inlineStr/text cell type) instead of silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.931.99'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.817.08'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''309.22'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.4654'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').Value = '''0.3654'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').Value = '''0.07367'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = '''0.8703'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').Value = '''20.25'
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('D12').Value = '1.805.69'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = '''5.387'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '''0.07111'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '''6.516'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').Value = '''91.25'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '''0.000008691'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '''14.64'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').Value = '26.956.28'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').Value = '''10.58'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('D24').Value = '2.059.01'
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').Value = '''1.895'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '''151.08'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').Value = '''18.40'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = '''2.138'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').Value = '''5.255'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').Value = '''116.34'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').Value = '''0.08880'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').Value = '''0.7591'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '''1.165'
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').Value = '''4.480'
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('D35').Value = '''2.897'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '''1.096'
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').Value = '''0.05290'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').Value = '''0.01947'
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('D40').Value = '''2.979'
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('D41').Value = '''0.5292'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').Value = '''7.157'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').Value = '''2.330'
$ws.Range('E43').Value = '  -3.82%  '
$ws.Range('D44').Value = '''0.1656'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = '''8.432'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('D46').Value = '''0.4851'
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').Value = '''10.41'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '''103.26'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').Value = '''1.661'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('D51').Value = '''0.06294'
$ws.Range('E51').Value = '  +0.02%  '
